# Daily update of covid19 tracker data files
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Country Updates")

# Bump "Updated on" date (column B) for all data rows (5-96) from 11-Apr-2020 to 12-Apr-2020
$ws.Range("B5:B96").Value = 43933

# Lithuania (row 34) narrative updates
# D34: Travel bans/restrictions
$ws.Range("D34").Value = ' - Foreigners are not allowed to enter the country. Exceptions include crews employed in Lithuanian companies engaged in international commercial transport or engaged in international commercial transport by any means of transport (unless they are found to have COVID-19 symptoms), foreigners who are immediate family members (parents, children, spouses, foster parents) to Lithuanian citizens, persons entitled to reside in Lithuania, diplomatic staff, NATO personnel.                                                                                                                                                            - Cruise ship access to Klaipeda seaport is prohibited.                                                                                                                                                                                      - From April 4 2020,  it is prohibited to transport passengers and/or passenger cars to Lithuania by ferries, except for the ferry on the route Kiel (Germany) – Klaipėda (Lithuania) – Kiel (Germany).  In exceptional cases,  the decision on the permit shall be made by the government.                                                                                                                                                       - From 4 April 2020 flights for the carriage of passengers to, and from international airports of Lithuania, shall be performed only upon obtaining a permit.                                                                     - Lithuanian consulates worldwide suspended acceptance and processing of visa applications, including visa applications in representation of another Schengen states
- Transit through the territory of Lithuania of persons returning to their country of residence and of crew members performing international commercial transport activities in all types of vehicles without the required stop is permitted. 
- A temporary ban on intercity travel over the Easter weekend was introduced. People will not be able to enter municipalities other than their place of residence (with some exceptions).                                                                                    The above restrictions are in effect between 16 and 27 April 2020. '

# E34: Closure of schools/universities
$ws.Range("E34").Value = ' - All educational institutions  (kindergartens, schools, universities, colleges, non-formal education institutions, etc.) closed from 13 March (as a recommendation) or from 16 March (compulsory suspension) until 13th April. The closure could be extended for a longer period.                                                                                              - On 30 March distance learning was introduced in Lithuania. The period of distance learning in schools and , vocational education and higher education institutions depends on the quarantine period. Distance learning is recommended for preschool and pre-primary education institutions but some kindergartens can operate and provide services for parents/families in need (i.e., doctors, police officers), upon municipality  decision. '

# F34: Cancellation of public events / Closure of public places
$ws.Range("F34").Value = '-All the events are prohibited throughout Lithuania. 
- Being in public spaces limited to 2 persons (except for members of immediate family).
- All shops are closed (with the exemption of supermarkets and pharmacies.'

# H34: Fiscal measures (overview)
$ws.Range("H34").Value = 'The government introduced on 16 March 2020 a fiscal package of EUR 2.5 billion (including EUR 0.5 billion for health) under the action plan to mitigate the impact of COVID-19. The measures aim to:                                                                                                                                                                                                        - Ensure the effective functioning of the health care system: the funds are dedicated, in particular, to the acquisition of equipment, finance additional health care costs, including bonuses to health care workers, and ensure subsidiary social guarantees for employees infected by a contagious disease.                                                                                                                                                                                                                                                                                                                                                                                                    - Help safeguard jobs and incomes, including subsidies wages, benefits for self-employed and funds to ensure the payment of sickness benefits of carers.                                                                                                                           -  Help business to save liquidity, including postponement of the payments date of tax arrears for the affected taxpayers.                                                                                                                                                                                 -Stimulate the economy including through  the acceleration of investment programmes (EUR 1.2 billion); reallocation of the EU investment funds to health, employment and business areas;  acceleration  of the use of the State budget appropriations for current expenditure; permission to use all funds of the Climate Change Programme, Road Maintenance and Development Programme, and to speed up the multi-apartment building renovation programme; and by offering possibilities for municipalities to implement the investment projects by extending lending targets and to ensure the emergency liquidity assistance.
                                                                                                                                                                                    To ensure the liquidity of the State Treasury, the action plan grants the right to the government to borrow additional EUR 5 billion (entered into force on 18 March).
                                                                                                                                                                                                                                                      As approved on 8 April, the government will provide short-term loans to municipalities and reimburse emergency management costs which were a result of central government decisions. Municipalities will also be able to have higher budget expenditures than earlier approved – EUR 30 million for small municipalities and EUR 10 million for large municipalities.'

# I34: Fiscal measures (income support)
$ws.Range("I34").Value = 'The fiscal package under the action plan to mitigate the impact of COVID-19, includes the following measures:                  
- Subsidising wages:The state will contribute to maintaining jobs with subsidies.. During downtime, employees  are paid at least the minimum wage, but only if the employment contract stipulates a full working time rate.  Subsidies‘ vary between 70% an  90% from the employee''s accrued salary, up to  a specified  limit.  Employers benefiting from the measure are obliged to keep at least 50% of the jobs in their firm for at least 3 months after the end of subsidy payments.
- Ensuring the payment of sickness benefits. it will be issued for up to 60 days, but no longer than the end of quarantine or emergency;                                                                                                                                                                          -  Increased sickness benefit for employees who became  infected during the performance of their duties. 
-  Flat rate sickness benefit (EUR 257) to self-employed persons who have  paid social insurance contributions for at least 3 months during the year prior to the declaration.
-  Scholarship for unemployed (0.39 of minimum monthly salary) even though the education process itself is temporarily inaccessible, or is carried out remotely, when the unemployed registered with the VET program. 
- Prolongation  (from 3 to 6 months) the grace period for repayments on mortgage loans to redundant persons.                                                                                                                                                                                              - Possibility of deferring or arranging in instalments the electricity and gas payments to public provider.                                                                                                                                                          - Recommendation for municipalities to exempt the taxpayers from taxes levied on business property and land. 
'

# J34: Fiscal measures (business support)
$ws.Range("J34").Value = 'The fiscal package under the action plan to mitigate the impact of COVID-19 includes the following measures:
 -Postponement of the payments date of tax arrears for the affected tax payers.
-Exemption from fines and default interest for failure to comply with tax obligations on time.                                                                                                  -Postponement of submission (and payment of) personal income tax returns and advanced corporate income tax returns.                              
 - Relief from import duties and VAT exemption on importation granted for goods needed to combat the effects of the COVID-19 outbreak pursuant to Commission’s decision.                                                                                                                                                                                    - The guarantee limit for the Agricultural Credit Guarantee Fund and ceiling for INVEGA (promoting the development of small- and medium-sized enterprises) guarantee are set to increase by EUR 500 million (entered into force on 18 March).
 - Possibility for business customers to defer or arrange the payments in instalments to the public provider of electricity and gas.
 - Recommendation for the municipalities to offer the possibility of deferring or arranging the payment of public utility charges and payments for heat in instalments. 
 - Recommendation for the municipalities to exempt the taxpayers from taxes levied on business property and land.                                                                                                                                     -  Allocation of EUR 1.3 billion to the firms facing liquidity and financial problems: EUR 851 million for state guaranties; EUR 287 million for loans; EUR 145 million for risk capital investment;  EUR 23 million for compensations of loans and leasing contracts interests (during the payment holiday up to 6 months); EUR 100 million for a Business Support Fund, and EUR 50 million for payable invoices loans.
'

